$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column G
$ws.Range("G1").Value = "avg_SQFT"

# Add AVERAGE formula for G2:G31
$ws.Range("G2:G31").Formula = "=AVERAGE(B2:C2)"

